$d = $word.ActiveDocument

# The header row of the first table contains "Month" and "Savings" cells
# whose run formatting needs to be refreshed (Bold/Italic/Strikethrough
# explicitly set), matching the POI 4.1.0 -> 5.2.3 upgrade fix.
$table = $d.Tables.Item(1)

for ($rowIdx = 1; $rowIdx -le $table.Rows.Count; $rowIdx++) {
    # Iterate the row's actual Cells collection (not Table.Cell(row,col)
    # across Columns.Count) so horizontally merged rows with fewer real
    # cells than Columns.Count are not addressed out of bounds.
    $row = $table.Rows.Item($rowIdx)
    for ($colIdx = 1; $colIdx -le $row.Cells.Count; $colIdx++) {
        $cell = $row.Cells.Item($colIdx)
        # Cell.Range.Text carries a trailing CR (13) + cell-mark (7); strip
        # those control characters before comparing the visible text.
        $cellText = $cell.Range.Text.TrimEnd([char]13, [char]7)

        if ($cellText -eq "Month" -or $cellText -eq "Savings") {
            # Restrict the range to the cell's actual content, excluding
            # the trailing cell-mark/paragraph-mark character, so only the
            # run's rPr is touched (not the paragraph mark's rPr).
            $contentRange = $cell.Range.Duplicate
            $contentRange.End = $contentRange.End - 1

            $contentRange.Font.Bold = $true
            $contentRange.Font.Italic = $false
            $contentRange.Font.Strikethrough = $false
        }
    }
}
